# PCA & GPU for new sim, Stairs, and hydrant adjustments
#
# This script reproduces, via the Excel COM object model, the edits that were
# made to data_manual.xlsx:
#   - Stairs positions (Connect_LDL_Front / Offload_Front) Location Y moved
#     from 756 to 825.
#   - Remove_GPU row's recorded location (I28/J28) cleared - GPU is no longer
#     parked at a fixed spot for the new sim.
#   - Attach_Tug row's recorded location/hydrant offsets (I33/J33/K33/L33,
#     including the K33 "960+(960-I33)" hydrant-offset formula) cleared.
#   - The special accent/theme font colour that had been applied to
#     alternating rows in column A is removed so the whole "Operation" column
#     reads in the normal/automatic text colour.
#   - The window/selection state is updated to where the author left off
#     (scrolled down to row 19, with I27 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Data edits
# ---------------------------------------------------------------------------

# Stairs (Connect_LDL_Front / Offload_Front): Location Y 756 -> 825
$ws.Range("J16").Value = 825
$ws.Range("J24").Value = 825

# Remove_GPU: clear the recorded Location X / Location Y
$ws.Range("I28").ClearContents()
$ws.Range("J28").ClearContents()

# Attach_Tug: clear the recorded Location X / Location Y / Location X2
# (hydrant offset formula) / Location Y2
$ws.Range("I33").ClearContents()
$ws.Range("J33").ClearContents()
$ws.Range("K33").ClearContents()
$ws.Range("L33").ClearContents()

# ---------------------------------------------------------------------------
# 2. Formatting - drop the accent/theme font colour used on alternating rows
#    of column A so the whole column uses the same (automatic/no colour)
#    formatting.
# ---------------------------------------------------------------------------
$ws.Range("A1:A36").Font.Color = 0

# ---------------------------------------------------------------------------
# 3. View / selection state
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)

# Window placement/size as last saved by the author
$win.Left = -28920
$win.Top = 2490
$win.Width = 29040
$win.Height = 15840

# Scroll the sheet so row 19 is at the top and select I27, matching where
# the author left the selection.
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("I27").Select()
